$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("I29").Value = 2.63
$ws.Range("L29").Value = 3.6
$ws.Range("W29").Value = 1.73
$ws.Range("X29").Value = 2.08
$ws.Range("Z29").Value = 1.47

# Row 30
$ws.Range("M30").Value = 1.14
$ws.Range("N30").Value = 5.5
$ws.Range("W30").Value = 1.67
$ws.Range("Z30").Value = 1.54

# Row 57
$ws.Range("AK57").Value = 1000

# Row 142
$ws.Range("S142").Value = 1.93
$ws.Range("T142").Value = 1.93

# Row 143
$ws.Range("Q143").Value = 1.88
$ws.Range("R143").Value = 1.98

# Row 146
$ws.Range("G146").Value = 1.65
$ws.Range("H146").Value = 3.7
$ws.Range("I146").Value = 4.75
$ws.Range("J146").Value = 2.18
$ws.Range("K146").Value = 2.2
$ws.Range("L146").Value = 4.85
$ws.Range("O146").Value = 1.23
$ws.Range("P146").Value = 3.35
$ws.Range("Q146").Value = 1.7
$ws.Range("R146").Value = 1.93
$ws.Range("U146").Value = 2.62
$ws.Range("V146").Value = 1.37
$ws.Range("AA146").Value = 7.5
$ws.Range("AD146").Value = 13
$ws.Range("AE146").Value = 12.5
$ws.Range("AF146").Value = 23
$ws.Range("AG146").Value = 11.5
$ws.Range("AH146").Value = 7.3
$ws.Range("AI146").Value = 15
$ws.Range("AL146").Value = 14
$ws.Range("AM146").Value = 29
$ws.Range("AN146").Value = 15
$ws.Range("AO146").Value = 90

# Row 147
$ws.Range("H147").Value = 4
$ws.Range("I147").Value = 5.8
$ws.Range("J147").Value = 1.98
$ws.Range("K147").Value = 2.32
$ws.Range("L147").Value = 5.4
$ws.Range("O147").Value = 1.2
$ws.Range("P147").Value = 3.6
$ws.Range("Q147").Value = 1.62
$ws.Range("R147").Value = 2.05
$ws.Range("U147").Value = 2.45
$ws.Range("AA147").Value = 7.6
$ws.Range("AC147").Value = 8
$ws.Range("AD147").Value = 10.75
$ws.Range("AG147").Value = 13
$ws.Range("AH147").Value = 8
$ws.Range("AL147").Value = 18
$ws.Range("AP147").Value = 55

# Row 230
$ws.Range("G230").Value = 1.75
$ws.Range("H230").Value = 3.9
$ws.Range("I230").Value = 4.5
$ws.Range("J230").Value = 2.3
$ws.Range("K230").Value = 2.25
$ws.Range("L230").Value = 4.75
$ws.Range("N230").Value = 12
$ws.Range("Q230").Value = 1.83
$ws.Range("R230").Value = 2.03
$ws.Range("W230").Value = 1.36
$ws.Range("X230").Value = 3
$ws.Range("Y230").Value = 1.8
$ws.Range("Z230").Value = 1.95
$ws.Range("AA230").Value = 7.5
$ws.Range("AG230").Value = 12

# Row 231
$ws.Range("S231").Value = 1.93
$ws.Range("T231").Value = 1.93

# Row 232
$ws.Range("G232").Value = 1.65
$ws.Range("H232").Value = 4.2
$ws.Range("I232").Value = 4.75
$ws.Range("J232").Value = 2.2
$ws.Range("K232").Value = 2.3
$ws.Range("L232").Value = 5
$ws.Range("O232").Value = 1.22
$ws.Range("P232").Value = 4
$ws.Range("Q232").Value = 1.75
$ws.Range("R232").Value = 2.05
$ws.Range("U232").Value = 2.75
$ws.Range("V232").Value = 1.4
$ws.Range("W232").Value = 1.33
$ws.Range("X232").Value = 3.25
$ws.Range("AB232").Value = 8
$ws.Range("AG232").Value = 13
$ws.Range("AH232").Value = 8
$ws.Range("AI232").Value = 17
$ws.Range("AM232").Value = 26
$ws.Range("AP232").Value = 41

# Row 233
$ws.Range("M233").Value = 1.04
$ws.Range("N233").Value = 13
$ws.Range("Q233").Value = 1.85
$ws.Range("R233").Value = 2
$ws.Range("U233").Value = 3
$ws.Range("V233").Value = 1.36

# Row 234
$ws.Range("G234").Value = 1.36
$ws.Range("H234").Value = 5
$ws.Range("I234").Value = 8
$ws.Range("J234").Value = 1.83
$ws.Range("K234").Value = 2.6
$ws.Range("L234").Value = 7
$ws.Range("N234").Value = 17
$ws.Range("O234").Value = 1.17
$ws.Range("P234").Value = 5
$ws.Range("Q234").Value = 1.57
$ws.Range("R234").Value = 2.35
$ws.Range("S234").Value = 1.93
$ws.Range("T234").Value = 1.93
$ws.Range("U234").Value = 2.38
$ws.Range("V234").Value = 1.53
$ws.Range("AB234").Value = 7
$ws.Range("AD234").Value = 9
$ws.Range("AH234").Value = 9.5
$ws.Range("AI234").Value = 19
$ws.Range("AK234").Value = 251
$ws.Range("AL234").Value = 21
$ws.Range("AP234").Value = 51
$ws.Range("AQ234").Value = 51

# Row 240
$ws.Range("G240").Value = 1.83
$ws.Range("H240").Value = 3.4
$ws.Range("I240").Value = 3.8
$ws.Range("J240").Value = 2.47
$ws.Range("K240").Value = 2.15
$ws.Range("W240").Value = 1.39
$ws.Range("X240").Value = 2.77
$ws.Range("Y240").Value = 1.7
$ws.Range("AA240").Value = 7.5
$ws.Range("AB240").Value = 9
$ws.Range("AE240").Value = 14.5
$ws.Range("AH240").Value = 6.8
$ws.Range("AI240").Value = 13.5
$ws.Range("AK240").Value = 400
$ws.Range("AL240").Value = 12.5
$ws.Range("AM240").Value = 23
$ws.Range("AO240").Value = 60
$ws.Range("AP240").Value = 32
$ws.Range("AQ240").Value = 35

# Row 241
$ws.Range("G241").Value = 2.65
$ws.Range("H241").Value = 3.15
$ws.Range("I241").Value = 2.47
$ws.Range("J241").Value = 3.25
$ws.Range("L241").Value = 3.15
$ws.Range("N241").Value = 7.1
$ws.Range("O241").Value = 1.31
$ws.Range("V241").Value = 1.31
$ws.Range("AC241").Value = 9.75
$ws.Range("AF241").Value = 29
$ws.Range("AG241").Value = 7.1
$ws.Range("AJ241").Value = 55
$ws.Range("AN241").Value = 9.5
$ws.Range("AO241").Value = 27
$ws.Range("AP241").Value = 21

# Row 242
$ws.Range("G242").Value = 1.72
$ws.Range("H242").Value = 3.6
$ws.Range("I242").Value = 4.1
$ws.Range("J242").Value = 2.32
$ws.Range("K242").Value = 2.18
$ws.Range("L242").Value = 4.55
$ws.Range("Q242").Value = 1.83
$ws.Range("R242").Value = 1.87
$ws.Range("U242").Value = 3
$ws.Range("V242").Value = 1.34
$ws.Range("W242").Value = 1.39
$ws.Range("X242").Value = 2.77
$ws.Range("Y242").Value = 1.8
$ws.Range("Z242").Value = 1.91
$ws.Range("AA242").Value = 7.1
$ws.Range("AB242").Value = 8.25
$ws.Range("AD242").Value = 13.5
$ws.Range("AH242").Value = 7.1
$ws.Range("AI242").Value = 15.5
$ws.Range("AJ242").Value = 70
$ws.Range("AL242").Value = 12
$ws.Range("AM242").Value = 23
$ws.Range("AN242").Value = 14
$ws.Range("AO242").Value = 65
$ws.Range("AP242").Value = 40
